$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 36.89194233333333
$ws.Range("N2").Value = 110.675827
$ws.Range("O2").Value = 0.3567095043190808
$ws.Range("P2").Value = 0.3567095043190809
$ws.Range("Q2").Value = 1.121662614702667
$ws.Range("R2").Value = 10.094963532324
$ws.Range("S2").Value = 0.3567095043190808
$ws.Range("T2").Value = 0.3567095043190809

$ws.Range("M3").Value = 42.68037399999999
$ws.Range("O3").Value = 0.4126780562577495
$ws.Range("P3").Value = 0.4126780562577496
$ws.Range("S3").Value = 0.4126780562577495
$ws.Range("T3").Value = 0.4126780562577496

$ws.Range("M4").Value = 23.85061433333334
$ws.Range("N4").Value = 71.55184300000001
$ws.Range("O4").Value = 0.2306124394231696
$ws.Range("P4").Value = 0.2306124394231696
$ws.Range("Q4").Value = 0.7251540781906668
$ws.Range("R4").Value = 6.526386703716001
$ws.Range("S4").Value = 0.2306124394231696
$ws.Range("T4").Value = 0.2306124394231696
